$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change -------------------------------------------------
# Clear the "DL01 Shortage" label from C3. This also causes the now-unused
# shared string to be dropped from the shared-strings table on save.
$ws.Range("C3").ClearContents()

# --- Column width changes -------------------------------------------
# Column D (Part Number) widened and no longer auto "best fit".
$ws.Columns("D").ColumnWidth = 17.5
# New explicit widths for columns G and H.
$ws.Columns("G").ColumnWidth = 9.833333333333334
$ws.Columns("H").ColumnWidth = 14.333333333333334

# --- View / selection changes ----------------------------------------
$win = $excel.ActiveWindow
# Reset the frozen/scrolled top-left cell back to A1 (removes topLeftCell="B1").
$win.ScrollColumn = 1
$win.ScrollRow = 1
# Zoom to 150%.
$win.Zoom = 150
# Move the selection to H1.
$ws.Range("H1").Select()
